# "Remove elective patients instead of CABG" — updates the SOFA flow-chart
# numbers on the single slide, plus two unused/orphaned locale date
# placeholders left over in the notes master and an Italian slide layout.

# Helper: replace a paragraph's whole text with a brand-new single run.
# A direct "$range.Text = $new" does a minimal character-level diff against
# the old text and keeps any matching prefix/suffix as separate runs; going
# through a throwaway value first avoids that and collapses the paragraph
# back down to one clean run (matching how PowerPoint behaves when a user
# selects the whole line and retypes it).
function Set-ParagraphText {
    param($range, [string]$newText)
    $range.Text = "~"
    $range.Text = $newText
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "73,181 adult ICU stays in MIMIC-IV" box: CABG -> Elective Admission,
#     and the GCS-missing count.
$excl1 = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParagraphText $excl1.Paragraphs(4) "Elective Admission (n=6,939)"
Set-ParagraphText $excl1.Paragraphs(5) "Information on GCS missing (n=1,146)"

# --- "xx,xxx adult ICU stays in MIMIC-IV / within first 24 hours" box.
$mimic24 = $s.Shapes.Item(5).TextFrame.TextRange
Set-ParagraphText $mimic24.Paragraphs(1) "33,968 adult ICU stays in MIMIC-IV "

# --- "200,859 ICU stays in eICU-CRD" exclusions box: CABG/sex-missing swap.
$excl2 = $s.Shapes.Item(7).TextFrame.TextRange
Set-ParagraphText $excl2.Paragraphs(4) "Information on sex missing and age < 18 (n=327)"
Set-ParagraphText $excl2.Paragraphs(5) "Elective Admission (n=19,384)"

# --- "127,380 adult ICU stays in eICU-CRD within first 24 hours" box.
$eicu24 = $s.Shapes.Item(11).TextFrame.TextRange
Set-ParagraphText $eicu24.Paragraphs(1) "113,118 adult ICU stays in eICU-CRD within first 24 hours"

# --- "LOS < 7 days (n=31,983)" box: only the digits after "(" get retyped,
#     so PowerPoint leaves the literal "(" in the first run and puts the
#     rest in a second run.
$los1 = $s.Shapes.Item(12).TextFrame.TextRange
$full = $los1.Text
$idx = $full.IndexOf("n=31,983)")
$sub = $los1.Characters($idx + 1, "n=31,983)".Length)
$sub.Text = "n=29,042)"

# --- "LOS < 7 days (n=112,810)" box.
$los2 = $s.Shapes.Item(16).TextFrame.TextRange
Set-ParagraphText $los2.Paragraphs(2) "LOS < 7 days (n=102,465)"

# --- "14,570 adult ICU stays in eICU-CRD / within first 7 days" box: only
#     the leading run's text changes, the rest of the paragraph (eICU-CRD,
#     line break, "within first 7 days") is untouched.
$eicu7 = $s.Shapes.Item(19).TextFrame.TextRange
$lead = $eicu7.Characters(1, "14,570 adult ICU stays in ".Length)
$lead.Text = "10,653 adult ICU stays in "

# --- Orphaned locale date placeholders (not shown on the actual slide, but
#     still present in the deck's notes master / Italian "Solo titolo"
#     layout) get bumped by two days along with everything else.
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "11.02.2023"

$itLayout = $p.SlideMaster.CustomLayouts.Item(16)
$itLayout.Shapes.Item(2).TextFrame.TextRange.Text = "11/02/2023"
